$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.711.54'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.680.77'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3912'
$ws.Range('E7').Value = '  -1.96%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3950'
$ws.Range('E8').Value = '  -2.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.003'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '51.81'
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.399'
$ws.Range('E11').Value = '  -5.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08654'
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.26'
$ws.Range('E13').Value = '  -3.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.324'
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.756'
$ws.Range('E15').Value = '  -4.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001314'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '1.645.98'
$ws.Range('E17').Value = '  -2.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.51'
$ws.Range('E18').Value = '  -3.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07083'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.25'
$ws.Range('E20').Value = '  -4.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.059'
$ws.Range('E21').Value = '  -2.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.004'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.91'
$ws.Range('E23').Value = '  -3.25%  '
$ws.Range('D24').Value = '24.696.62'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.347'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.833'
$ws.Range('E26').Value = '  -3.33%  '
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.62'
$ws.Range('E28').Value = '  -2.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.837'
$ws.Range('E29').Value = '  -6.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '145.83'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.869'
$ws.Range('E31').Value = '  -5.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.378'
$ws.Range('E32').Value = '  +6.58%  '
$ws.Range('D33').Value = '1.864.07'
$ws.Range('E33').Value = '  -2.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08405'
$ws.Range('E34').Value = '  -4.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.03050'
$ws.Range('E35').Value = '  -5.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.956'
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9993'
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2793'
$ws.Range('E38').Value = '  -2.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09437'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.60'
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.542'
$ws.Range('E41').Value = '  +5.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7917'
$ws.Range('E42').Value = '  -7.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.52'
$ws.Range('E43').Value = '  -4.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.53'
$ws.Range('E44').Value = '  -5.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7143'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('E46').Value = '  -5.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.180'
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.08666'
$ws.Range('E48').Value = '  +3.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.340'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '137.80'
$ws.Range('E51').Value = '  -2.22%  '
